$wb = $excel.ActiveWorkbook

# ---- DataSet sheet: replace old QA admin credentials with new Admin (Priya) credentials ----
$ds = $wb.Worksheets.Item("DataSet")
$ds.Range("B2").Value = "pjampala@helenoftroy.com"
$ds.Range("C2").Value = "pjampala@helenoftroy.com"
$ds.Range("D2").Value = "Priya@123"
$ds.Range("E2").Value = "Priya@123"
$ds.Range("H2").Value = "pjampala@helenoftroy.com"
$ds.Range("I2").Value = "pjampala@helenoftroy.com"
$ds.Range("B3").Value = "pjampala@helenoftroy.com"
$ds.Range("D3").Value = "Priya@123"
$ds.Range("E3").Value = "Priya@123"
$ds.Range("H3").Value = "pjampala@helenoftroy.com"
$ds.Range("I3").Value = "pjampala@helenoftroy.com"
$ds.Range("H17").Value = "pjampala@helenoftroy.com"
$ds.Range("D18").Value = "Priya@123"
$ds.Range("E18").Value = "Priya@1234"
$ds.Range("H20").Value = "pjampala@helenoftroy.com"
$ds.Range("L26").Value = "OXO 5-Piece Barware Set - Exclusive"
$ds.Range("H33").Value = "pjampala@helenoftroy.com"
$ds.Range("H34").Value = "pjampala@helenoftroy.com"
$ds.Range("I34").Value = "pjampala@helenoftroy.com"
$ds.Range("H35").Value = "pjampala@helenoftroy.com"
$ds.Range("I35").Value = "pjampala@helenoftroy.com"
$ds.Range("B39").Value = "pjampala@helenoftroy.com"
$ds.Range("C39").Value = "pjampala@helenoftroy.com"
$ds.Range("D39").Value = "Priya@123"
$ds.Range("E39").Value = "Priya@123"
$ds.Range("H39").Value = "pjampala@helenoftroy.com"
$ds.Range("I39").Value = "pjampala@helenoftroy.com"
$ds.Range("B40").Value = "pjampala@helenoftroy.com"
$ds.Range("C40").Value = "pjampala@helenoftroy.com"
$ds.Range("D40").Value = "Priya@123"
$ds.Range("E40").Value = "Priya@123"
$ds.Range("H40").Value = "pjampala@helenoftroy.com"
$ds.Range("I40").Value = "pjampala@helenoftroy.com"
$ds.Range("B41").Value = "pjampala@helenoftroy.com"
$ds.Range("C41").Value = "pjampala@helenoftroy.com"
$ds.Range("D41").Value = "Priya@123"
$ds.Range("E41").Value = "Priya@123"
$ds.Range("H41").Value = "pjampala@helenoftroy.com"
$ds.Range("I41").Value = "pjampala@helenoftroy.com"
$ds.Range("B42").Value = "pjampala@helenoftroy.com"
$ds.Range("C42").Value = "pjampala@helenoftroy.com"
$ds.Range("D42").Value = "Priya@123"
$ds.Range("E42").Value = "Priya@123"
$ds.Range("H42").Value = "pjampala@helenoftroy.com"
$ds.Range("I42").Value = "pjampala@helenoftroy.com"

# ---- E2E sheet: replace old QA admin credentials with new Admin (Priya) credentials ----
$e2e = $wb.Worksheets.Item("E2E")
$e2e.Range("B2").Value = "pjampala@helenoftroy.com"
$e2e.Range("D2").Value = "Priya@123"
$e2e.Range("E2").Value = "Priya@123"
$e2e.Range("H2").Value = "pjampala@helenoftroy.com"
$e2e.Range("D18").Value = "Priya@123"
$e2e.Range("E18").Value = "Priya@1234"
$e2e.Range("H34").Value = "pjampala@helenoftroy.com"
$e2e.Range("H35").Value = "pjampala@helenoftroy.com"
$e2e.Range("B39").Value = "pjampala@helenoftroy.com"
$e2e.Range("D39").Value = "Priya@123"
$e2e.Range("E39").Value = "Priya@123"
$e2e.Range("H39").Value = "pjampala@helenoftroy.com"
$e2e.Range("B40").Value = "pjampala@helenoftroy.com"
$e2e.Range("D40").Value = "Priya@123"
$e2e.Range("E40").Value = "Priya@123"
$e2e.Range("H40").Value = "pjampala@helenoftroy.com"
$e2e.Range("B41").Value = "pjampala@helenoftroy.com"
$e2e.Range("D41").Value = "Priya@123"
$e2e.Range("E41").Value = "Priya@123"
$e2e.Range("H41").Value = "pjampala@helenoftroy.com"
$e2e.Range("B42").Value = "pjampala@helenoftroy.com"
$e2e.Range("D42").Value = "Priya@123"
$e2e.Range("E42").Value = "Priya@123"
$e2e.Range("H42").Value = "pjampala@helenoftroy.com"
$e2e.Range("B61").Value = "pjampala@helenoftroy.com"
$e2e.Range("D61").Value = "Priya@123"
$e2e.Range("E61").Value = "Priya@123"
$e2e.Range("H61").Value = "pjampala@helenoftroy.com"

# ---- E2E sheet: new "preprodURL" column header (copy header style from neighboring DOB header) ----
$e2e.Range("AJ1").Copy()
$e2e.Range("AK1").PasteSpecial(-4122)
$e2e.Range("AK1").Value = "preprodURL"

# ---- E2E sheet: new row 63 - Admin login details for E2E validation ----
$e2e.Range("A63").Value = "Login Details"
$e2e.Range("B2").Copy()
$e2e.Range("B63").PasteSpecial(-4122)
$e2e.Range("B63").Value = "pjampala@helenoftroy.com"
$e2e.Range("D63").Value = "Acmshnqpimqk9("
$e2e.Range("D2").Copy()
$e2e.Range("AK63").PasteSpecial(-4122)
$e2e.Range("AK63").Value = "https://na-preprod.hele.digital/heledigitaladmin/admin/"

# ---- New hyperlinks for the Admin login details row ----
$e2e.Hyperlinks.Add($e2e.Range("AK63"), "https://na-preprod.hele.digital/heledigitaladmin/admin/") | Out-Null
$e2e.Hyperlinks.Add($e2e.Range("B63"), "mailto:pjampala@helenoftroy.com") | Out-Null

# ---- Make DataSet the active sheet (matches final saved view state) ----
$ds.Activate()
